$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A6").Value = "Напиток"
$ws.Range("B6").Value = "Кофе"
$ws.Range("C6").Value = "Американо или эспрессо на выбор, 0.2л"
$ws.Range("D6").Value = 150
$ws.Range("E6").Value = "https://www.pngmart.com/files/21/Coffee-Cup-PNG-Isolated-Pic.png"

$ws.Range("A7").Value = "Напиток"
$ws.Range("C7").Value = "Зеленый или черный на выбор"
$ws.Range("B7").Value = "Чай Lipton"
$ws.Range("D7").Value = 125
$ws.Range("E7").Value = "https://nnjfood.ru/upload/iblock/9dc/ce0wwpjev5mcg1qzbst7v562jrbs6cs3.jpg"

$ws.Range("G18").Select()
